# Scheduled market-data refresh for the Leve profit tables.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# with freshly pulled Universalis price data, row by row, per job class sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40 (Horn Glue)
$ws.Range("H40").Value = 1558.9445
$ws.Range("I40").Value = 1551.2222
$ws.Range("J40").Value = 1566.6666
$ws.Range("K40").Value = 1551.2222
$ws.Range("L40").Value = 1566.6666
$ws.Range("M40").Value = -1376.2222
$ws.Range("N40").Value = -1916.6666

# Row 70 (Holy Water)
$ws.Range("H70").Value = 60159.117
$ws.Range("I70").Value = 251075.5
$ws.Range("J70").Value = 1415.6154
$ws.Range("K70").Value = 753226.5
$ws.Range("L70").Value = 4246.8462
$ws.Range("M70").Value = -752956.5
$ws.Range("N70").Value = -4786.8462

# Row 73 (Holy Water)
$ws.Range("H73").Value = 60159.117
$ws.Range("I73").Value = 251075.5
$ws.Range("J73").Value = 1415.6154
$ws.Range("K73").Value = 753226.5
$ws.Range("L73").Value = 4246.8462
$ws.Range("M73").Value = -752290.5
$ws.Range("N73").Value = -6118.8462

# Row 129 (Commanding Craftsman's Draught)
$ws.Range("H129").Value = 1209.0513
$ws.Range("J129").Value = 1434.2
$ws.Range("L129").Value = 4302.6
$ws.Range("N129").Value = -14302.6

# Row 137 (Magnesia Whetstone)
$ws.Range("H137").Value = 1176.5151
$ws.Range("I137").Value = 1013.68
$ws.Range("J137").Value = 1685.375
$ws.Range("K137").Value = 3041.04
$ws.Range("L137").Value = 5056.125
$ws.Range("M137").Value = -491.04
$ws.Range("N137").Value = -10156.125

# Row 138 (Cunning Craftsman's Tisane)
$ws.Range("H138").Value = 2247.4927
$ws.Range("I138").Value = 1261.6123
$ws.Range("J138").Value = 4662.9
$ws.Range("K138").Value = 3784.8369
$ws.Range("L138").Value = 13988.7
$ws.Range("M138").Value = 1355.1631
$ws.Range("N138").Value = -24268.7

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Steel Ingot)
$ws.Range("H32").Value = 9953.477000000001
$ws.Range("I32").Value = 9859.966
$ws.Range("J32").Value = 10873
$ws.Range("K32").Value = 9859.966
$ws.Range("L32").Value = 10873
$ws.Range("M32").Value = -9572.966
$ws.Range("N32").Value = -11447

# Row 123 (High Durium Armguards of Maiming)
$ws.Range("H123").Value = 24281.75
$ws.Range("J123").Value = 24281.75
$ws.Range("L123").Value = 24281.75
$ws.Range("N123").Value = -34081.75

# Row 131 (Chondrite Top of Maiming)
$ws.Range("H131").Value = 47673
$ws.Range("J131").Value = 47673
$ws.Range("L131").Value = 47673
$ws.Range("N131").Value = -57753

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (Iron Ingot)
$ws.Range("H20").Value = 20989.744
$ws.Range("I20").Value = 26515.59
$ws.Range("K20").Value = 26515.59
$ws.Range("M20").Value = -26268.59

# Row 45 (Wrapped Steel Awl)
$ws.Range("H45").Value = 22074
$ws.Range("J45").Value = 22074
$ws.Range("L45").Value = 22074
$ws.Range("N45").Value = -23690

# Row 86 (Adamantite Nugget)
$ws.Range("H86").Value = 254051.12
$ws.Range("I86").Value = 5920.4
$ws.Range("J86").Value = 667602.3
$ws.Range("K86").Value = 5920.4
$ws.Range("L86").Value = 667602.3
$ws.Range("M86").Value = -4797.4
$ws.Range("N86").Value = -669848.3

# Row 89 (Adamantite Nugget)
$ws.Range("H89").Value = 254051.12
$ws.Range("I89").Value = 5920.4
$ws.Range("J89").Value = 667602.3
$ws.Range("K89").Value = 29602
$ws.Range("L89").Value = 3338011.5
$ws.Range("M89").Value = -23986
$ws.Range("N89").Value = -3349243.5

# Row 134 (Ruthenium Ingot)
$ws.Range("H134").Value = 2534.4614
$ws.Range("I134").Value = 2263.4375
$ws.Range("J134").Value = 3773.4285
$ws.Range("K134").Value = 6790.3125
$ws.Range("L134").Value = 11320.2855
$ws.Range("M134").Value = -4255.3125
$ws.Range("N134").Value = -16390.2855

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Walnut Lumber)
$ws.Range("H31").Value = 2346.6943
$ws.Range("I31").Value = 1338.3043
$ws.Range("K31").Value = 1338.3043
$ws.Range("M31").Value = -1043.3043

# Row 34 (Walnut Lumber)
$ws.Range("H34").Value = 2346.6943
$ws.Range("I34").Value = 1338.3043
$ws.Range("K34").Value = 1338.3043
$ws.Range("M34").Value = -1136.3043

# Row 37 (Yew Crook)
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 92 (Gyr Abanian Flour)
$ws.Range("H92").Value = 583
$ws.Range("I92").Value = 666
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 1998
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = -750
$ws.Range("N92").Value = -3996

# Row 133 (Boiled Alpaca Steak)
$ws.Range("H133").Value = 5619.522
$ws.Range("I133").Value = 1202.625
$ws.Range("J133").Value = 7975.2
$ws.Range("K133").Value = 3607.875
$ws.Range("L133").Value = 23925.6
$ws.Range("M133").Value = 1452.125
$ws.Range("N133").Value = -34045.6

# Row 138 (Tacos Al Pastor)
$ws.Range("H138").Value = 2007.9412
$ws.Range("I138").Value = 920.8182
$ws.Range("J138").Value = 4001
$ws.Range("K138").Value = 2762.4546
$ws.Range("L138").Value = 12003
$ws.Range("M138").Value = 2377.5454
$ws.Range("N138").Value = -22283

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (Hardsilver Ingot)
$ws.Range("H80").Value = 3055.5557
$ws.Range("I80").Value = 2928.5715
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 2928.5715
$ws.Range("L80").Value = 3500
$ws.Range("M80").Value = -1930.5715
$ws.Range("N80").Value = -5496

# Row 83 (Hardsilver Ingot)
$ws.Range("H83").Value = 3055.5557
$ws.Range("I83").Value = 2928.5715
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 14642.8575
$ws.Range("L83").Value = 17500
$ws.Range("M83").Value = -9650.8575
$ws.Range("N83").Value = -27484

# Row 109 (Hematite Earrings of Healing)
$ws.Range("H109").Value = 9049.333000000001
$ws.Range("J109").Value = 9049.333000000001
$ws.Range("L109").Value = 9049.333000000001
$ws.Range("N109").Value = -11129.333

# Row 123 (Ametrine Ring of Fending)
$ws.Range("H123").Value = 10516.5
$ws.Range("J123").Value = 10516.5
$ws.Range("L123").Value = 10516.5
$ws.Range("N123").Value = -15416.5

$ws = $wb.Worksheets.Item("LTW")
# Row 93 (Gagana Leather)
$ws.Range("H93").Value = 1191.8
$ws.Range("I93").Value = 1026.6428
$ws.Range("K93").Value = 1026.6428
$ws.Range("M93").Value = 221.3571999999999

# Row 139 (Gomphotherium Doublet of Gathering)
$ws.Range("H139").Value = 44524.168
$ws.Range("J139").Value = 44524.168
$ws.Range("L139").Value = 44524.168
$ws.Range("N139").Value = -54804.168

$ws = $wb.Worksheets.Item("WVR")
# Row 62 (Rainbow Cloth)
$ws.Range("H62").Value = 4428.5713

# Row 65 (Rainbow Cloth)
$ws.Range("H65").Value = 4428.5713

# Row 123 (Fingerless Darkhempen Gloves of Healing)
$ws.Range("H123").Value = 24636.334
$ws.Range("J123").Value = 24636.334
$ws.Range("L123").Value = 24636.334
$ws.Range("N123").Value = -34436.334

# Row 126 (Snow Linen)
$ws.Range("H126").Value = 7593.75
$ws.Range("I126").Value = 8193.182000000001
$ws.Range("K126").Value = 24579.546
$ws.Range("M126").Value = -22109.546

# Row 139 (Rroneek Serge Trousers of Gathering)
$ws.Range("H139").Value = 69645
$ws.Range("J139").Value = 69645
$ws.Range("L139").Value = 69645
$ws.Range("N139").Value = -79925
